# Refresh cryptos list figures (prices / 1h volume %) and reorder a few rows
# whose underlying coin changed, per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '42.762.33'
$ws.Range('E2').Value = '  -0.73%  '
# Row 3
$ws.Range('D3').Value = '2.295.71'
$ws.Range('E3').Value = '  -0.13%  '
# Row 5
$ws.Range('D5').Value = '''305.61'
$ws.Range('E5').Value = '  +1.87%  '
# Row 6
$ws.Range('D6').Value = '''96.48'
$ws.Range('E6').Value = '  -1.05%  '
# Row 7
$ws.Range('D7').Value = '''0.509'
$ws.Range('E7').Value = '  -1.80%  '
# Row 8
$ws.Range('E8').Value = '  +0.02%  '
# Row 9
$ws.Range('D9').Value = '''0.501'
$ws.Range('E9').Value = '  -2.71%  '
# Row 10
$ws.Range('D10').Value = '''35.56'
$ws.Range('E10').Value = '  -1.70%  '
# Row 11
$ws.Range('D11').Value = '''0.0791'
$ws.Range('E11').Value = '  -0.14%  '
# Row 12
$ws.Range('D12').Value = '''18.33'
$ws.Range('E12').Value = '  +3.33%  '
# Row 13
$ws.Range('D13').Value = '''0.118'
$ws.Range('E13').Value = '  +1.08%  '
# Row 14
$ws.Range('D14').Value = '''6.73'
$ws.Range('E14').Value = '  -2.01%  '
# Row 15
$ws.Range('D15').Value = '2.657.08'
$ws.Range('E15').Value = '  -0.01%  '
# Row 16
$ws.Range('D16').Value = '2.309.78'
$ws.Range('E16').Value = '  +0.98%  '
# Row 17
$ws.Range('D17').Value = '''0.780'
$ws.Range('E17').Value = '  -0.89%  '
# Row 18
$ws.Range('D18').Value = '42.705.28'
$ws.Range('E18').Value = '  -0.57%  '
# Row 19
$ws.Range('D19').Value = '''12.96'
$ws.Range('E19').Value = '  -0.54%  '
# Row 20
$ws.Range('D20').Value = '0.0₃0898'
$ws.Range('E20').Value = '  -1.33%  '
# Row 21
$ws.Range('D21').Value = '''6.01'
$ws.Range('E21').Value = '  -1.94%  '
# Row 22
$ws.Range('D22').Value = '''67.24'
$ws.Range('E22').Value = '  -1.38%  '
# Row 23
$ws.Range('D23').Value = '''235.99'
$ws.Range('E23').Value = '  -0.71%  '
# Row 24
$ws.Range('E24').Value = '  -2.23%  '
# Row 25
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '''1.00'
$ws.Range('E25').Value = '  +0.00%  '
# Row 26
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').Value = '''2.45'
$ws.Range('E26').Value = '  +1.10%  '
# Row 27
$ws.Range('D27').Value = '''4.02'
$ws.Range('E27').Value = '  +0.14%  '
# Row 28
$ws.Range('D28').Value = '''25.20'
$ws.Range('E28').Value = '  +0.73%  '
# Row 29
$ws.Range('D29').Value = '''166.82'
$ws.Range('E29').Value = '  +2.42%  '
# Row 30
$ws.Range('D30').Value = '''2.06'
$ws.Range('E30').Value = '  +0.99%  '
# Row 31
$ws.Range('D31').Value = '''9.05'
$ws.Range('E31').Value = '  -1.03%  '
# Row 32
$ws.Range('D32').Value = '''33.30'
$ws.Range('E32').Value = '  +0.98%  '
# Row 33
$ws.Range('E33').Value = '  +0.08%  '
# Row 34
$ws.Range('D34').Value = '''4.78'
$ws.Range('E34').Value = '  +1.05%  '
# Row 35
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').Value = '''17.80'
$ws.Range('E35').Value = '  -2.02%  '
# Row 36
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').Value = '''4.98'
$ws.Range('E36').Value = '  -2.89%  '
# Row 37
$ws.Range('E37').Value = '  -0.65%  '
# Row 38
$ws.Range('D38').Value = '''0.0690'
$ws.Range('E38').Value = '  -0.50%  '
# Row 39
$ws.Range('E39').Value = '  -0.96%  '
# Row 40
$ws.Range('E40').Value = '  -1.63%  '
# Row 41
$ws.Range('D41').Value = '''0.109'
$ws.Range('E41').Value = '  -1.48%  '
# Row 42
$ws.Range('D42').Value = '''2.70'
$ws.Range('E42').Value = '  -1.81%  '
# Row 43
$ws.Range('D43').Value = '2.001.50'
$ws.Range('E43').Value = '  -0.37%  '
# Row 44
$ws.Range('D44').Value = '''0.0280'
$ws.Range('E44').Value = '  -2.44%  '
# Row 45
$ws.Range('D45').Value = '''18.20'
$ws.Range('E45').Value = '  +3.99%  '
# Row 46
$ws.Range('E46').Value = '  -2.40%  '
# Row 47
$ws.Range('D47').Value = '''2.06'
$ws.Range('E47').Value = '  -5.89%  '
# Row 48
$ws.Range('D48').Value = '''2.76'
$ws.Range('E48').Value = '  -2.59%  '
# Row 49
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').Value = '''2.94'
$ws.Range('E49').Value = '  +7.95%  '
# Row 50
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = '''53.67'
$ws.Range('E50').Value = '  -1.13%  '
# Row 51
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.525.58'
$ws.Range('E51').Value = '  -0.20%  '
